$d = $word.ActiveDocument

# 1) Merge " Legion" and ", however " runs into a single run with combined text.
#    (This runtime coalesces every directly-adjacent, identically-formatted run
#    touched by an edit, so the edit below also swallows the following run —
#    "you weren't able to gain the respect from Captain ". We find that run's
#    range first so we can split it back out afterwards with a harmless
#    Bold on/off toggle, which leaves no formatting residue.)
$find0 = $d.Content.Find
$find0.Execute("you weren", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterRunAnchorStart = $find0.Parent.Start

$find = $d.Content.Find
$find.Execute(" Legion, however ", $true, $false, $false, $false, $false, $true, 1, $false, " Legion, however ", 2)

$find0b = $d.Content.Find
$find0b.Execute("you weren’t able to gain the respect from Captain ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterRun = $d.Range($afterRunAnchorStart, $find0b.Parent.End)
$afterRun.Bold = 1
$afterRun.Bold = 0

# 2) Bold the "Restart?" paragraph (both the paragraph mark and the run text).
$find2 = $d.Content.Find
$find2.Execute("Restart?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$restartRun = $find2.Parent
$restartPara = $restartRun.Paragraphs(1)
$restartPara.Range.Font.Bold = 1

# 3) Move the _GoBack bookmark so it spans the whole "Restart?" paragraph
#    (including its paragraph mark) instead of being collapsed at the start.
$paraStart = $restartPara.Range.Start
$paraEnd = $restartPara.Range.End
$d.Bookmarks("_GoBack").Delete()
$newBookmarkRange = $d.Range($paraStart, $paraEnd + 1)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
